$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(3)

# Insert a new "property_category" column before the existing "date" column
# (H), pushing date/legislator_name/legislator_id one column to the right.
$ws.Range("H1").EntireColumn.Insert()

$ws.Range("H1").Value = "property_category"
$ws.Range("H2").Value = "stock"
$ws.Range("H3").Value = "stock"

# Clean up the stock name (drop the spaces in "Equinox Minerals Limited").
$ws.Range("B2").Value = "★EquinoxMineralsLimited"
